$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "(409, 233)"

$ws.Range("C2").Value = "(419, 250)"
$ws.Range("E2").Value = 11.40175425099138
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 9

$ws.Range("C3").Value = "(414, 239)"
$ws.Range("D3").Value = $true
$ws.Range("E3").Value = 7.810249675906654
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 6
